$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data: insert two new rows at the top of the data block (rows 2-3),
# pushing all existing rows down by two. Shared columns (A,B,C,E-K,Q,R,T) stay
# identical across every row in this sheet, so copy them from the (now-shifted)
# row 4 which still holds the template values.
$ws.Rows("2:3").Insert(-4121)

# The insert copies formatting down from the header row (bold) - strip that,
# then restore the date number format on column D to match the rest of the
# "Fecha" column.
$ws.Range("A2:T3").ClearFormats()
$ws.Range("D2:D3").NumberFormat = $ws.Range("D4").NumberFormat

# Row 2: new "Primera" quality entry for 2022-02-15 (serial 44607)
$ws.Range("A2").Value = $ws.Range("A4").Value2
$ws.Range("B2").Value = $ws.Range("B4").Value2
$ws.Range("C2").Value = $ws.Range("C4").Value2
$ws.Range("D2").Value = 44607
$ws.Range("E2").Value = $ws.Range("E4").Value2
$ws.Range("F2").Value = $ws.Range("F4").Value2
$ws.Range("G2").Value = $ws.Range("G4").Value2
$ws.Range("H2").Value = $ws.Range("H4").Value2
$ws.Range("I2").Value = $ws.Range("I4").Value2
$ws.Range("J2").Value = $ws.Range("J4").Value2
$ws.Range("K2").Value = $ws.Range("K4").Value2
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11500
$ws.Range("Q2").Value = $ws.Range("Q4").Value2
$ws.Range("R2").Value = $ws.Range("R4").Value2
$ws.Range("S2").Value = 639
$ws.Range("T2").Value = $ws.Range("T4").Value2

# Row 3: new "Segunda" quality entry for 2022-02-15 (serial 44607)
$ws.Range("A3").Value = $ws.Range("A4").Value2
$ws.Range("B3").Value = $ws.Range("B4").Value2
$ws.Range("C3").Value = $ws.Range("C4").Value2
$ws.Range("D3").Value = 44607
$ws.Range("E3").Value = $ws.Range("E4").Value2
$ws.Range("F3").Value = $ws.Range("F4").Value2
$ws.Range("G3").Value = $ws.Range("G4").Value2
$ws.Range("H3").Value = $ws.Range("H4").Value2
$ws.Range("I3").Value = $ws.Range("I4").Value2
$ws.Range("J3").Value = $ws.Range("J4").Value2
$ws.Range("K3").Value = $ws.Range("K4").Value2
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 240
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 9500
$ws.Range("Q3").Value = $ws.Range("Q4").Value2
$ws.Range("R3").Value = $ws.Range("R4").Value2
$ws.Range("S3").Value = 528
$ws.Range("T3").Value = $ws.Range("T4").Value2

Write-Host "done"
